$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters: it controls the order new shared strings get appended to the
# shared string table, which must match the target workbook's table order.
$ws.Range("E5").Value = "Radioisotopos"
$ws.Range("A43").Value = "Gabriel Barros G de Souza"
$ws.Range("A23").Value = "Helena Nery Alves Pinto"
$ws.Range("A10").Value = "Juliana Leal"
$ws.Range("A35").Value = "Marcos Paulo Maia Jorge"
$ws.Range("A44").Value = "Mariana M Vale"
$ws.Range("A13").Value = "Rhuanna Cavalcante Paulo"
$ws.Range("C13").Value = "Externo"
$ws.Range("E26").Value = "NA"
$ws.Range("E45").Value = "Radioisotopos"

$ws.Range("E27").Select()
